# FCA_RYG.xlsx dashboard edit
# - Clean up stray trailing non-breaking-space characters from the
#   "P_PMSC_PT_IntgAbsExhAccel" parameter name in D3.
# - Bold the header row (row 1).
# - Leave the active selection on D20 (matches the author's last click
#   before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the D3 cell value: same text, but without the trailing
#    non-breaking spaces / space that had accumulated on the string.
$ws.Range("D3").Value = "P_PMSC_PT_IntgAbsExhAccel"

# 2. Bold the header row (A1:H1)
$ws.Range("A1:H1").Font.Bold = $true

# 3. Update the saved selection/active cell
$ws.Range("D20").Select() | Out-Null
